$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 - new task "18. Change colour of underlined words"
$ws.Range("A19").Value = "18. Change colour of underlined words"
$ws.Range("B19").Value = "Hight"
$ws.Range("C19").Value = "Done"
$ws.Range("C19").Font.Color = 5287936
$ws.Range("D19").Value = "Arthur"

# Row 20 - new task "19. Make feature to add word to ban list by clicking"
$ws.Range("A20").Value = "19. Make feature to add word to ban list by clicking"
$ws.Range("B20").Value = "Middle"
$ws.Range("C20").Value = "Open"

# Move active selection to A20, matching the saved sheet view state
$ws.Range("A20").Select()
